$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (Covered_States), shifting the
# existing E:AC columns one position to the right (-> F:AD). Excel's
# native "insert column" behaviour copies formatting (incl. width) from
# the column to the left, which is what the target file shows (new E
# picks up D's bold/centered header style and plain data style).
$ws.Columns("E").Insert()

# New column E mirrors column D's width as closely as this host allows.
$ws.Columns("E").ColumnWidth = 23.25

# Populate the new header/value pair.
$ws.Range("E1").Value = "Covered_States"
$ws.Range("E2").Value = "{vendor:covered_state}"

# Restore a normal (top-left) view and select F12, matching the saved
# workbook view/selection in the edited file.
$ws.Range("F12").Select()

# Set the page to portrait orientation (Page Layout change recorded in
# the edited file).
$ws.PageSetup.Orientation = 1
